$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values in column D and E are written as literal text
# (matching the workbook's existing inline-string / text format), not converted to numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "303.18"
$ws.Range("E2").Value = "0.11%"

# Row 3
$ws.Range("D3").Value = "37.15"
$ws.Range("E3").Value = "6.62%"

# Row 4
$ws.Range("D4").Value = "4.995"
$ws.Range("E4").Value = "-3.44%"

# Row 5
$ws.Range("D5").Value = "0.07825"
$ws.Range("E5").Value = "0.82%"

# Row 6
$ws.Range("D6").Value = "2.201"
$ws.Range("E6").Value = "-3.12%"

# Row 7
$ws.Range("D7").Value = "8.014"
$ws.Range("E7").Value = "0.06%"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9158"
$ws.Range("E8").Value = "-1.25%"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.09680"
$ws.Range("E9").Value = "-4.20%"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1880"
$ws.Range("E10").Value = "3.58%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08581"
$ws.Range("E11").Value = "0.27%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03567"
$ws.Range("E12").Value = "3.02%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09968"
$ws.Range("E13").Value = "0.69%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001484"
$ws.Range("E14").Value = "0.15%"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005669"
$ws.Range("E15").Value = "-2.15%"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.458"
$ws.Range("E16").Value = "-0.33%"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.046"
$ws.Range("E17").Value = "1.36%"

# Row 18
$ws.Range("D18").Value = "2.376"
$ws.Range("E18").Value = "12.92%"

# Row 19
$ws.Range("E19").Value = "0.66%"

# Row 20
$ws.Range("D20").Value = "0.1311"
$ws.Range("E20").Value = "-1.19%"

# Row 21
$ws.Range("D21").Value = "4.781"
$ws.Range("E21").Value = "4.83%"

# Row 22
$ws.Range("E22").Value = "-1.94%"

# Row 23
$ws.Range("D23").Value = "0.04617"
$ws.Range("E23").Value = "0.02%"

# Row 24
$ws.Range("D24").Value = "0.001234"
$ws.Range("E24").Value = "0.91%"

# Row 25
$ws.Range("D25").Value = "0.004783"
$ws.Range("E25").Value = "8.07%"

# Row 26
$ws.Range("D26").Value = "0.0001402"
$ws.Range("E26").Value = "7.60%"

# Row 27
$ws.Range("E27").Value = "38.94%"

# Row 39
$ws.Range("D39").Value = "0.01771"
$ws.Range("E39").Value = "0.98%"

# Row 40
$ws.Range("D40").Value = "0.04743"
$ws.Range("E40").Value = "0.63%"

# Row 41
$ws.Range("D41").Value = "0.008076"
$ws.Range("E41").Value = "6.35%"

# Row 42
$ws.Range("E42").Value = "-0.78%"

# Row 43
$ws.Range("D43").Value = "0.007831"
$ws.Range("E43").Value = "13.91%"

# Row 44
$ws.Range("D44").Value = "0.002093"
$ws.Range("E44").Value = "-5.55%"

# Row 45
$ws.Range("D45").Value = "0.009960"
$ws.Range("E45").Value = "8.42%"

# Row 46
$ws.Range("D46").Value = "0.00006117"
$ws.Range("E46").Value = "2.34%"

# Row 47
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.08%"

# Row 48
$ws.Range("D48").Value = "7.978"
$ws.Range("E48").Value = "191.97%"

# Row 49
$ws.Range("E49").Value = "-0.51%"

# Row 50
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").Value = "0.08%"

# Row 51
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").Value = "0.08%"

